$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 2.431118488311768
$ws.Range("B1").Value = 2.723590850830078
$ws.Range("C1").Value = 1.625187993049622
$ws.Range("D1").Value = 1.319416284561157
$ws.Range("E1").Value = 1.230550169944763
